$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 32: Automata for the People
$ws.Range("H32").Value = 7147564
$ws.Range("I32").Value = 7403
$ws.Range("J32").Value = 9094881
$ws.Range("K32").Value = 7403
$ws.Range("L32").Value = 9094881
$ws.Range("M32").Value = -7077
$ws.Range("N32").Value = -9095533

# Row 68: Can't Sleep, Inquisitors Will Eat Me
$ws.Range("H68").Value = 62038.5
$ws.Range("I68").Value = 24156
$ws.Range("J68").Value = 74666
$ws.Range("K68").Value = 24156
$ws.Range("L68").Value = 74666
$ws.Range("M68").Value = -23407
$ws.Range("N68").Value = -76164

# Row 71: Allow No Fallacies (L)
$ws.Range("H71").Value = 62038.5
$ws.Range("I71").Value = 24156
$ws.Range("J71").Value = 74666
$ws.Range("K71").Value = 72468
$ws.Range("L71").Value = 223998
$ws.Range("M71").Value = -68724
$ws.Range("N71").Value = -231486

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 1305.625
$ws.Range("J88").Value = 455.25
$ws.Range("L88").Value = 455.25
$ws.Range("N88").Value = -1267.25

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 1305.625
$ws.Range("J91").Value = 455.25
$ws.Range("L91").Value = 455.25
$ws.Range("N91").Value = -3263.25

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 296.5
$ws.Range("I107").Value = 294.75
$ws.Range("K107").Value = 294.75
$ws.Range("M107").Value = 1625.25

# Row 125: Body over Mind
$ws.Range("H125").Value = 1746.7
$ws.Range("I125").Value = 1567.25
$ws.Range("K125").Value = 14105.25
$ws.Range("M125").Value = -11645.25


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4920.148
$ws.Range("I61").Value = 3991.75
$ws.Range("J61").Value = 7572.7144
$ws.Range("K61").Value = 3991.75
$ws.Range("L61").Value = 7572.7144
$ws.Range("M61").Value = -3779.75
$ws.Range("N61").Value = -7996.7144

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 6668.3335
$ws.Range("I63").Value = 6502.5
$ws.Range("J63").Value = 7000
$ws.Range("K63").Value = 6502.5
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = -5816.5
$ws.Range("N63").Value = -8372

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 6668.3335
$ws.Range("I66").Value = 6502.5
$ws.Range("J66").Value = 7000
$ws.Range("K66").Value = 32512.5
$ws.Range("L66").Value = 35000
$ws.Range("M66").Value = -29080.5
$ws.Range("N66").Value = -41864

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4920.148
$ws.Range("I136").Value = 3991.75
$ws.Range("J136").Value = 7572.7144
$ws.Range("K136").Value = 11975.25
$ws.Range("L136").Value = 22718.1432
$ws.Range("M136").Value = -9425.25
$ws.Range("N136").Value = -27818.1432


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 6194
$ws.Range("I134").Value = 4696.5
$ws.Range("J134").Value = 7283.091
$ws.Range("K134").Value = 14089.5
$ws.Range("L134").Value = 21849.273
$ws.Range("M134").Value = -11554.5
$ws.Range("N134").Value = -26919.273


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 109: Playing the Market
$ws.Range("H109").Value = 40624.375
$ws.Range("J109").Value = 40624.375
$ws.Range("L109").Value = 40624.375
$ws.Range("N109").Value = -42704.375

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 4213761.5
$ws.Range("I122").Value = 9310209
$ws.Range("K122").Value = 27930627
$ws.Range("M122").Value = -27928177

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3493.9678
$ws.Range("I134").Value = 2242.25
$ws.Range("J134").Value = 7785.5713
$ws.Range("K134").Value = 6726.75
$ws.Range("L134").Value = 23356.7139
$ws.Range("M134").Value = -4191.75
$ws.Range("N134").Value = -28426.7139


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 6: Meat-lover's Special
$ws.Range("H6").Value = 1400.3334
$ws.Range("I6").Value = 1400.3334
$ws.Range("K6").Value = 4201.0002
$ws.Range("M6").Value = -4088.0002

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 1223.7693
$ws.Range("J113").Value = 1355.375
$ws.Range("L113").Value = 4066.125
$ws.Range("N113").Value = -8406.125


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 42823.645
$ws.Range("I80").Value = 72420.69
$ws.Range("J80").Value = 3360.9167
$ws.Range("K80").Value = 72420.69
$ws.Range("L80").Value = 3360.9167
$ws.Range("M80").Value = -71422.69
$ws.Range("N80").Value = -5356.9167

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 42823.645
$ws.Range("I83").Value = 72420.69
$ws.Range("J83").Value = 3360.9167
$ws.Range("K83").Value = 362103.45
$ws.Range("L83").Value = 16804.5835
$ws.Range("M83").Value = -357111.45
$ws.Range("N83").Value = -26788.5835

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 6667.2354
$ws.Range("I113").Value = 3668.182
$ws.Range("J113").Value = 12165.5
$ws.Range("K113").Value = 3668.182
$ws.Range("L113").Value = 12165.5
$ws.Range("M113").Value = -1498.182
$ws.Range("N113").Value = -16505.5

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 3333.8096
$ws.Range("I122").Value = 2263.25
$ws.Range("K122").Value = 6789.75
$ws.Range("M122").Value = -4339.75


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 6: Sticking Their Necks Out
$ws.Range("H6").Value = 36333
$ws.Range("J6").Value = 36333
$ws.Range("L6").Value = 36333
$ws.Range("N6").Value = -36557

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 4824.357
$ws.Range("J68").Value = 4835.143
$ws.Range("L68").Value = 4835.143
$ws.Range("N68").Value = -6333.143

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 4824.357
$ws.Range("J71").Value = 4835.143
$ws.Range("L71").Value = 24175.715
$ws.Range("N71").Value = -31663.715

# Row 102: Shrug It On
$ws.Range("H102").Value = 48999.2
$ws.Range("J102").Value = 48999.2
$ws.Range("L102").Value = 48999.2
$ws.Range("N102").Value = -55489.2

# Row 103: Security Breeches
$ws.Range("H103").Value = 27333
$ws.Range("J103").Value = 27333
$ws.Range("L103").Value = 27333
$ws.Range("N103").Value = -29677

# Row 104: Brace Yourselves
$ws.Range("H104").Value = 21285.715
$ws.Range("J104").Value = 21285.715
$ws.Range("L104").Value = 21285.715
$ws.Range("N104").Value = -28273.715


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 4374.5
$ws.Range("I62").Value = 3749
$ws.Range("K62").Value = 3749
$ws.Range("M62").Value = -3125

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 4374.5
$ws.Range("I65").Value = 3749
$ws.Range("K65").Value = 18745
$ws.Range("M65").Value = -15625

# Row 74: Clothing the Naked Truth
$ws.Range("H74").Value = 8654.4
$ws.Range("J74").Value = 8363
$ws.Range("L74").Value = 8363
$ws.Range("N74").Value = -10235

# Row 77: When in Robes (L)
$ws.Range("H77").Value = 8654.4
$ws.Range("J77").Value = 8363
$ws.Range("L77").Value = 25089
$ws.Range("N77").Value = -34449

